# The document has two logo pictures (the Pearson logo in the footers,
# the BTEC logo in the headers) that each appear twice: once in the
# "first page" header/footer and once in the "default" (other pages)
# header/footer. This edit simply swaps each picture's stored file name
# the other logo's/its own sibling's name, i.e.
#   footers: "image1.png" -> "image2.png"
#   headers: "image2.jpg" -> "image1.jpg"
#
# InlineShape (the object Word hands back for a picture that sits
# in-line with text, which is how these logos are inserted) has no
# writable Name property in Word's object model - only a floating
# Shape does. So for each logo we briefly convert the inline picture to
# a floating shape, rename it, then convert it straight back to an
# inline shape (restoring the original inline layout).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoPicture {
    param($range)

    $inlinePic = $range.InlineShapes.Item(1)

    # Shape.Name is settable (InlineShape.Name is not), so hop over to
    # the floating-shape form just long enough to rename, then convert
    # straight back to keep the original inline (wp:inline) layout.
    $floatingPic = $inlinePic.ConvertToShape()

    $oldName = $floatingPic.Name
    if ($oldName -eq "image1.png") {
        $newName = "image2.png"
    } elseif ($oldName -eq "image2.jpg") {
        $newName = "image1.jpg"
    } else {
        $newName = $oldName
    }
    $floatingPic.Name = $newName

    $floatingPic.ConvertToInlineShape() | Out-Null

    Write-Output ("Renamed picture '" + $oldName + "' -> '" + $newName + "'")
}

# Default (odd/other pages) footer + first-page footer - Pearson logo.
Rename-LogoPicture $sec.Footers.Item(1).Range
Rename-LogoPicture $sec.Footers.Item(2).Range

# Default (odd/other pages) header + first-page header - BTEC logo.
Rename-LogoPicture $sec.Headers.Item(1).Range
Rename-LogoPicture $sec.Headers.Item(2).Range
